$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (column C) date for every existing data row
#    (rows 2..394) from 2023-09-23 (45192) to 2023-10-03 (45202).
for ($r = 2; $r -le 394; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 394 previously lacked the explicit row-height marker that every
# other data row carries; restore it so the row again has an explicit
# customHeight of 15.
$ws.Rows.Item(394).RowHeight = 15

# 2. Append the two new cleared-notification rows.

# Row 395
$ws.Cells.Item(395, 1).Value = "A 46566-2023"
$ws.Cells.Item(395, 2).Value = 45198
$ws.Cells.Item(395, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(395, 3).Value = 45202
$ws.Cells.Item(395, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(395, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(395, 5).Value = "ÖSTRA GÖINGE"
$ws.Cells.Item(395, 7).Value = 4.2
$ws.Cells.Item(395, 8).Value = 0
$ws.Cells.Item(395, 9).Value = 0
$ws.Cells.Item(395, 10).Value = 0
$ws.Cells.Item(395, 11).Value = 0
$ws.Cells.Item(395, 12).Value = 0
$ws.Cells.Item(395, 13).Value = 0
$ws.Cells.Item(395, 14).Value = 0
$ws.Cells.Item(395, 15).Value = 0
$ws.Cells.Item(395, 16).Value = 0
$ws.Cells.Item(395, 17).Value = 0
$ws.Cells.Item(395, 18).WrapText = $true
$ws.Rows.Item(395).RowHeight = 15

# Row 396
$ws.Cells.Item(396, 1).Value = "A 46677-2023"
$ws.Cells.Item(396, 2).Value = 45198
$ws.Cells.Item(396, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(396, 3).Value = 45202
$ws.Cells.Item(396, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(396, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(396, 5).Value = "ÖSTRA GÖINGE"
$ws.Cells.Item(396, 7).Value = 0.5
$ws.Cells.Item(396, 8).Value = 0
$ws.Cells.Item(396, 9).Value = 0
$ws.Cells.Item(396, 10).Value = 0
$ws.Cells.Item(396, 11).Value = 0
$ws.Cells.Item(396, 12).Value = 0
$ws.Cells.Item(396, 13).Value = 0
$ws.Cells.Item(396, 14).Value = 0
$ws.Cells.Item(396, 15).Value = 0
$ws.Cells.Item(396, 16).Value = 0
$ws.Cells.Item(396, 17).Value = 0
$ws.Cells.Item(396, 18).WrapText = $true
